# Update "想去人数" (F column) values on sheet "展览" and "全部类型"
# to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 81
$ws1.Range("F6").Value  = 545
$ws1.Range("F7").Value  = 7683
$ws1.Range("F9").Value  = 204
$ws1.Range("F10").Value = 1080
$ws1.Range("F11").Value = 669
$ws1.Range("F12").Value = 15
$ws1.Range("F13").Value = 27
$ws1.Range("F14").Value = 175
$ws1.Range("F16").Value = 203
$ws1.Range("F17").Value = 760

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 81
$ws4.Range("F7").Value  = 545
$ws4.Range("F8").Value  = 7683
$ws4.Range("F10").Value = 204
$ws4.Range("F11").Value = 1080
$ws4.Range("F12").Value = 669
$ws4.Range("F13").Value = 15
$ws4.Range("F14").Value = 27
$ws4.Range("F15").Value = 175
$ws4.Range("F17").Value = 203
$ws4.Range("F18").Value = 760
